$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# (losing formatting such as trailing zeros, e.g. "27.00" -> 27).
$textCells = @(
    'D5',
    'D6',
    'D14',
    'D20',
    'D21',
    'D22',
    'D24',
    'D25',
    'D29',
    'D31',
    'D35',
    'D38',
    'D40',
    'D45',
    'D46',
    'D49',
    'D50'
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '57.712.12'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '3.100.57'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '523.52'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('D6').Value = '141.54'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.099.42'
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('E12').Value = '  +2.62%  '
$ws.Range('D13').Value = '3.632.80'
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').Value = '0.131'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = '57.785.81'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '3.098.93'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D20').Value = '12.78'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').Value = '8.04'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').Value = '338.90'
$ws.Range('E22').Value = '  +1.85%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '0.511'
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('D25').Value = '66.61'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '0.0₃0916'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = '6.48'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D31').Value = '7.17'
$ws.Range('E31').Value = '  -1.25%  '
$ws.Range('E32').Value = '  +2.67%  '
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').Value = '155.67'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('E36').Value = '  +1.44%  '
$ws.Range('E37').Value = '  +1.56%  '
$ws.Range('D38').Value = '27.00'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').Value = '0.0661'
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('D42').Value = '3.138.70'
$ws.Range('E42').Value = '  +0.98%  '
$ws.Range('E43').Value = '  +10.56%  '
$ws.Range('E44').Value = '  +3.66%  '
$ws.Range('D45').Value = '36.81'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').Value = '2.294.71'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('D49').Value = '0.985'
$ws.Range('E49').Value = '  +4.48%  '
$ws.Range('D50').Value = '20.47'
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('E51').Value = '  +1.35%  '
